$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing columns right
# (Player ID moves from A to B, etc.)
$ws.Columns("A:A").Insert()

# New column A is "Match ID": header in row 2, bold (no border) style applied
# to the header + all data rows, and a literal value of 1 in every data row.
$ws.Range("A2").Value = "Match ID"
$ws.Range("A2:A19").Font.Bold = $true
$ws.Range("A4:A19").Value = 1

# Row 20 (the hidden totals row) also gets a Match ID value, but keeps the
# default (unbolded) style. Writing directly into a hidden row stamps a
# bogus custom row height on this engine, so temporarily unhide it first.
$ws.Rows("20:20").Hidden = $false
$ws.Range("A20").Value = 1
$ws.Rows("20:20").Hidden = $true

# Selection / active cell moves as part of the edit.
$ws.Range("F25").Select() | Out-Null
